$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.772266626358032
$ws.Range("B1").Value = 2.601320505142212
$ws.Range("C1").Value = 2.89776086807251
$ws.Range("D1").Value = 3.912126302719116
$ws.Range("E1").Value = 4.896999835968018
